# EPBDS-10595: Improve @ServiceExtraMethod to use generated beans from the
# rules. Arrays support is added.
#
# Adds two new example blocks to the sheet demonstrating the overloaded
# "myRule" SmartRules method that now accepts an array parameter
# (MyDatatype[]) plus an extra String parameter, and the brand new
# "myRule2" method with the same array-based signature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: myRule(MyDatatype[] inputParam, String x) -> rows 25-29 ---
$ws.Range("C25").Value = "SmartRules Double myRule(MyDatatype[] inputParam, String x)"

$ws.Range("C26").Value = "x"
$ws.Range("D26").Value = "result"

$ws.Range("C27").Value = "a"
$ws.Range("D27").Value = 100

$ws.Range("C28").Value = "a"
$ws.Range("D28").Value = 200

$ws.Range("C29").Value = "c"
$ws.Range("D29").Value = 300

# --- Block 2: myRule2(MyDatatype[] inputParam, String x) -> rows 35-39 ---
$ws.Range("C35").Value = "SmartRules Double myRule2(MyDatatype[] inputParam, String x)"

$ws.Range("C36").Value = "x"
$ws.Range("D36").Value = "result"

$ws.Range("C37").Value = "a"
$ws.Range("D37").Value = 100

$ws.Range("C38").Value = "a"
$ws.Range("D38").Value = 200

$ws.Range("C39").Value = "c"
$ws.Range("D39").Value = 300

# Leave the selection where the author left it before saving.
$ws.Range("L13").Select()
